$p = $ppt.ActivePresentation

# --- 1. Swap the custom table style applied to the three tables that were
#        using the "Table_0" style ({F3E65DCF-...}) over to the built-in
#        style {EAF3956F-36E6-4F8B-A166-11E38433DE1F}. In each of slides
#        14, 15 and 16 the table is the first shape on the slide. -----------
$tableStyleId = "{EAF3956F-36E6-4F8B-A166-11E38433DE1F}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($tableStyleId)
    }
}

# --- 2. Swap the two themes in the deck: the design theme that slides /
#        the slide master use ("Integral" / "Red Violet") becomes the
#        plain "Office" palette. (The master/slides share a single live
#        theme color scheme, so editing it from any slide updates it
#        everywhere.) ---------------------------------------------------
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

# index -> (R,G,B) for the "Office" color scheme (dk1,lt1,dk2,lt2,
# accent1-6,hlink,folHlink)
$officeColors = @(
    @(0x00,0x00,0x00),  # 1  dk1
    @(0xFF,0xFF,0xFF),  # 2  lt1
    @(0x44,0x54,0x6A),  # 3  dk2
    @(0xE7,0xE6,0xE6),  # 4  lt2
    @(0x5B,0x9B,0xD5),  # 5  accent1
    @(0xED,0x7D,0x31),  # 6  accent2
    @(0xA5,0xA5,0xA5),  # 7  accent3
    @(0xFF,0xC0,0x00),  # 8  accent4
    @(0x44,0x72,0xC4),  # 9  accent5
    @(0x70,0xAD,0x47),  # 10 accent6
    @(0x05,0x63,0xC1),  # 11 hlink
    @(0x95,0x4F,0x72)   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $rgbBytes = $officeColors[$i - 1]
    $r = $rgbBytes[0]
    $g = $rgbBytes[1]
    $b = $rgbBytes[2]
    $val = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($i).RGB = $val
}
